# Refresh market-price derived columns (H-N) across the Sheets workbook
# to match the latest scheduled data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 649.7646999999999
$ws.Cells.Item(98, 9).Value = 649.7646999999999
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 649.7646999999999
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 848.2353000000001
$ws.Cells.Item(98, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 649.7646999999999
$ws.Cells.Item(122, 9).Value = 649.7646999999999
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 1949.2941
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = 500.7059000000002
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 12421723
$ws.Cells.Item(132, 9).Value = 14538422
$ws.Cells.Item(132, 10).Value = 250704
$ws.Cells.Item(132, 11).Value = 43615266
$ws.Cells.Item(132, 12).Value = 752112
$ws.Cells.Item(132, 13).Value = -43612736
$ws.Cells.Item(132, 14).Value = -757172

$ws.Cells.Item(138, 8).Value = 4661.375
$ws.Cells.Item(138, 9).Value = 1276.1111
$ws.Cells.Item(138, 10).Value = 5789.7964
$ws.Cells.Item(138, 11).Value = 3828.3333
$ws.Cells.Item(138, 12).Value = 17369.3892
$ws.Cells.Item(138, 13).Value = 1311.6667
$ws.Cells.Item(138, 14).Value = -27649.3892

$ws.Cells.Item(141, 8).Value = 55811.75
$ws.Cells.Item(141, 9).Value = 28263.055
$ws.Cells.Item(141, 11).Value = 84789.16500000001
$ws.Cells.Item(141, 13).Value = -79609.16500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2705.3403
$ws.Cells.Item(32, 9).Value = 2381.1333
$ws.Cells.Item(32, 10).Value = 10000
$ws.Cells.Item(32, 11).Value = 2381.1333
$ws.Cells.Item(32, 12).Value = 10000
$ws.Cells.Item(32, 13).Value = -2094.1333
$ws.Cells.Item(32, 14).Value = -10574

$ws.Cells.Item(61, 8).Value = 2710.8462
$ws.Cells.Item(61, 9).Value = 2556.6667
$ws.Cells.Item(61, 10).Value = 3558.8333
$ws.Cells.Item(61, 11).Value = 2556.6667
$ws.Cells.Item(61, 12).Value = 3558.8333
$ws.Cells.Item(61, 13).Value = -2344.6667
$ws.Cells.Item(61, 14).Value = -3982.8333

$ws.Cells.Item(132, 8).Value = 3223.611
$ws.Cells.Item(132, 9).Value = 2814.0938
$ws.Cells.Item(132, 11).Value = 8442.2814
$ws.Cells.Item(132, 13).Value = -5912.2814

$ws.Cells.Item(136, 8).Value = 2710.8462
$ws.Cells.Item(136, 9).Value = 2556.6667
$ws.Cells.Item(136, 10).Value = 3558.8333
$ws.Cells.Item(136, 11).Value = 7670.000100000001
$ws.Cells.Item(136, 12).Value = 10676.4999
$ws.Cells.Item(136, 13).Value = -5120.000100000001
$ws.Cells.Item(136, 14).Value = -15776.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 32170.334
$ws.Cells.Item(82, 9).Value = 32170.334
$ws.Cells.Item(82, 11).Value = 32170.334
$ws.Cells.Item(82, 13).Value = -31787.334

$ws.Cells.Item(85, 8).Value = 32170.334
$ws.Cells.Item(85, 9).Value = 32170.334
$ws.Cells.Item(85, 11).Value = 32170.334
$ws.Cells.Item(85, 13).Value = -30844.334

$ws.Cells.Item(134, 8).Value = 4933
$ws.Cells.Item(134, 9).Value = 4800
$ws.Cells.Item(134, 11).Value = 14400
$ws.Cells.Item(134, 13).Value = -11865

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 3350.6667
$ws.Cells.Item(2, 9).Value = 2019.8
$ws.Cells.Item(2, 10).Value = 10005
$ws.Cells.Item(2, 11).Value = 2019.8
$ws.Cells.Item(2, 12).Value = 10005
$ws.Cells.Item(2, 13).Value = -1906.8
$ws.Cells.Item(2, 14).Value = -10231

$ws.Cells.Item(31, 8).Value = 2149967.2
$ws.Cells.Item(31, 9).Value = 3167438
$ws.Cells.Item(31, 10).Value = 7923.8423
$ws.Cells.Item(31, 11).Value = 3167438
$ws.Cells.Item(31, 12).Value = 7923.8423
$ws.Cells.Item(31, 13).Value = -3167143
$ws.Cells.Item(31, 14).Value = -8513.8423

$ws.Cells.Item(34, 8).Value = 2149967.2
$ws.Cells.Item(34, 9).Value = 3167438
$ws.Cells.Item(34, 10).Value = 7923.8423
$ws.Cells.Item(34, 11).Value = 3167438
$ws.Cells.Item(34, 12).Value = 7923.8423
$ws.Cells.Item(34, 13).Value = -3167236
$ws.Cells.Item(34, 14).Value = -8327.8423

$ws.Cells.Item(58, 8).Value = 3449
$ws.Cells.Item(58, 9).Value = 2907.1
$ws.Cells.Item(58, 11).Value = 2907.1
$ws.Cells.Item(58, 13).Value = -2704.1

$ws.Cells.Item(60, 8).Value = 0
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 11).Value = 0
$ws.Cells.Item(60, 13).ClearContents()

$ws.Cells.Item(99, 8).Value = 3713.5386
$ws.Cells.Item(99, 9).Value = 3525.0908
$ws.Cells.Item(99, 10).Value = 4750
$ws.Cells.Item(99, 11).Value = 3525.0908
$ws.Cells.Item(99, 12).Value = 4750
$ws.Cells.Item(99, 13).Value = -2027.0908
$ws.Cells.Item(99, 14).Value = -7746

$ws.Cells.Item(126, 8).Value = 3713.5386
$ws.Cells.Item(126, 9).Value = 3525.0908
$ws.Cells.Item(126, 10).Value = 4750
$ws.Cells.Item(126, 11).Value = 10575.2724
$ws.Cells.Item(126, 12).Value = 14250
$ws.Cells.Item(126, 13).Value = -8105.2724
$ws.Cells.Item(126, 14).Value = -19190

$ws.Cells.Item(135, 8).Value = 162695
$ws.Cells.Item(135, 10).Value = 162695
$ws.Cells.Item(135, 12).Value = 162695
$ws.Cells.Item(135, 14).Value = -172835

$ws.Cells.Item(136, 8).Value = 3449
$ws.Cells.Item(136, 9).Value = 2907.1
$ws.Cells.Item(136, 11).Value = 8721.299999999999
$ws.Cells.Item(136, 13).Value = -6171.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 82.393936
$ws.Cells.Item(2, 9).Value = 72.90000000000001
$ws.Cells.Item(2, 11).Value = 437.4
$ws.Cells.Item(2, 13).Value = -324.4

$ws.Cells.Item(5, 8).Value = 6675.5
$ws.Cells.Item(5, 10).Value = 9106.857
$ws.Cells.Item(5, 12).Value = 27320.571
$ws.Cells.Item(5, 14).Value = -27544.571

$ws.Cells.Item(68, 8).Value = 1788990.6
$ws.Cells.Item(68, 9).Value = 2829
$ws.Cells.Item(68, 10).Value = 2503455.2
$ws.Cells.Item(68, 11).Value = 8487
$ws.Cells.Item(68, 12).Value = 7510365.600000001
$ws.Cells.Item(68, 13).Value = -7676
$ws.Cells.Item(68, 14).Value = -7511987.600000001

$ws.Cells.Item(71, 8).Value = 1788990.6
$ws.Cells.Item(71, 9).Value = 2829
$ws.Cells.Item(71, 10).Value = 2503455.2
$ws.Cells.Item(71, 11).Value = 25461
$ws.Cells.Item(71, 12).Value = 22531096.8
$ws.Cells.Item(71, 13).Value = -21405
$ws.Cells.Item(71, 14).Value = -22539208.8

$ws.Cells.Item(113, 8).Value = 339.07693
$ws.Cells.Item(113, 9).Value = 359.41177
$ws.Cells.Item(113, 10).Value = 300.66666
$ws.Cells.Item(113, 11).Value = 1078.23531
$ws.Cells.Item(113, 12).Value = 901.9999799999999
$ws.Cells.Item(113, 13).Value = 1091.76469
$ws.Cells.Item(113, 14).Value = -5241.99998

$ws.Cells.Item(121, 8).Value = 1727.4073
$ws.Cells.Item(121, 9).Value = 484.7
$ws.Cells.Item(121, 10).Value = 5278
$ws.Cells.Item(121, 11).Value = 1454.1
$ws.Cells.Item(121, 12).Value = 15834
$ws.Cells.Item(121, 13).Value = -144.0999999999999
$ws.Cells.Item(121, 14).Value = -18454

$ws.Cells.Item(124, 8).Value = 19153.918
$ws.Cells.Item(124, 10).Value = 20148.428
$ws.Cells.Item(124, 12).Value = 60445.284
$ws.Cells.Item(124, 14).Value = -70265.284

$ws.Cells.Item(129, 8).Value = 7617798.5
$ws.Cells.Item(129, 9).Value = 16500759
$ws.Cells.Item(129, 10).Value = 3832.1428
$ws.Cells.Item(129, 11).Value = 49502277
$ws.Cells.Item(129, 12).Value = 11496.4284
$ws.Cells.Item(129, 13).Value = -49497277
$ws.Cells.Item(129, 14).Value = -21496.4284

$ws.Cells.Item(131, 8).Value = 24851.445
$ws.Cells.Item(131, 10).Value = 2921.3057
$ws.Cells.Item(131, 12).Value = 8763.917099999999
$ws.Cells.Item(131, 14).Value = -18843.9171

$ws.Cells.Item(132, 8).Value = 1520.1818
$ws.Cells.Item(132, 10).Value = 1978.5714
$ws.Cells.Item(132, 12).Value = 17807.1426
$ws.Cells.Item(132, 14).Value = -22867.1426

$ws.Cells.Item(133, 8).Value = 11166.5

$ws.Cells.Item(135, 8).Value = 6675.5
$ws.Cells.Item(135, 10).Value = 9106.857
$ws.Cells.Item(135, 12).Value = 81961.713
$ws.Cells.Item(135, 14).Value = -87031.713

$ws.Cells.Item(137, 8).Value = 5337.636
$ws.Cells.Item(137, 9).Value = 5412.6665
$ws.Cells.Item(137, 10).Value = 5000
$ws.Cells.Item(137, 11).Value = 16237.9995
$ws.Cells.Item(137, 12).Value = 15000
$ws.Cells.Item(137, 13).Value = -11137.9995
$ws.Cells.Item(137, 14).Value = -25200

$ws.Cells.Item(140, 8).Value = 2350.7036
$ws.Cells.Item(140, 9).Value = 2268.0386
$ws.Cells.Item(140, 10).Value = 4500
$ws.Cells.Item(140, 11).Value = 6804.1158
$ws.Cells.Item(140, 12).Value = 13500
$ws.Cells.Item(140, 13).Value = -1624.1158
$ws.Cells.Item(140, 14).Value = -23860

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 9449.357
$ws.Cells.Item(7, 9).Value = 10024.25
$ws.Cells.Item(7, 11).Value = 10024.25
$ws.Cells.Item(7, 13).Value = -9912.25

$ws.Cells.Item(40, 8).Value = 5909.909
$ws.Cells.Item(40, 9).Value = 5334.3335
$ws.Cells.Item(40, 11).Value = 5334.3335
$ws.Cells.Item(40, 13).Value = -5198.3335

$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 13).ClearContents()

$ws.Cells.Item(104, 8).Value = 37653.6
$ws.Cells.Item(104, 10).Value = 37653.6
$ws.Cells.Item(104, 12).Value = 37653.6
$ws.Cells.Item(104, 14).Value = -44641.6

$ws.Cells.Item(105, 8).Value = 40000
$ws.Cells.Item(105, 9).Value = 40000
$ws.Cells.Item(105, 11).Value = 40000
$ws.Cells.Item(105, 13).Value = -36506

$ws.Cells.Item(126, 8).Value = 9449.357
$ws.Cells.Item(126, 9).Value = 10024.25
$ws.Cells.Item(126, 11).Value = 30072.75
$ws.Cells.Item(126, 13).Value = -27602.75

$ws.Cells.Item(132, 8).Value = 3992.7585
$ws.Cells.Item(132, 9).Value = 3415
$ws.Cells.Item(132, 10).Value = 9000
$ws.Cells.Item(132, 11).Value = 10245
$ws.Cells.Item(132, 12).Value = 27000
$ws.Cells.Item(132, 13).Value = -7715
$ws.Cells.Item(132, 14).Value = -32060

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 19524
$ws.Cells.Item(62, 9).Value = 18823.646
$ws.Cells.Item(62, 10).Value = 22500.5
$ws.Cells.Item(62, 11).Value = 18823.646
$ws.Cells.Item(62, 12).Value = 22500.5
$ws.Cells.Item(62, 13).Value = -18199.646
$ws.Cells.Item(62, 14).Value = -23748.5

$ws.Cells.Item(65, 8).Value = 19524
$ws.Cells.Item(65, 9).Value = 18823.646
$ws.Cells.Item(65, 10).Value = 22500.5
$ws.Cells.Item(65, 11).Value = 94118.23000000001
$ws.Cells.Item(65, 12).Value = 112502.5
$ws.Cells.Item(65, 13).Value = -90998.23000000001
$ws.Cells.Item(65, 14).Value = -118742.5

$ws.Cells.Item(81, 8).Value = 11571.286
$ws.Cells.Item(81, 9).Value = 35833
$ws.Cells.Item(81, 10).Value = 4954.4546
$ws.Cells.Item(81, 11).Value = 71666
$ws.Cells.Item(81, 12).Value = 9908.9092
$ws.Cells.Item(81, 13).Value = -70605
$ws.Cells.Item(81, 14).Value = -12030.9092

$ws.Cells.Item(84, 8).Value = 11571.286
$ws.Cells.Item(84, 9).Value = 35833
$ws.Cells.Item(84, 10).Value = 4954.4546
$ws.Cells.Item(84, 11).Value = 358330
$ws.Cells.Item(84, 12).Value = 49544.546
$ws.Cells.Item(84, 13).Value = -353026
$ws.Cells.Item(84, 14).Value = -60152.546

$ws.Cells.Item(132, 8).Value = 2185.1538
$ws.Cells.Item(132, 9).Value = 2185.1538
$ws.Cells.Item(132, 11).Value = 6555.4614
$ws.Cells.Item(132, 13).Value = -4025.4614
